# Updates cryptos list: refresh Price/Volume(1h) columns for each coin row,
# plus two swapped-rank pairs (rows 28/29 and 47/48) whose Coin/Link/Price/
# Volume cells moved to the other row of the pair.
# Price values ("D" column) that look like plain numbers are entered with a
# leading apostrophe so Excel stores them as literal text (matching the
# original inlineStr cell contents, e.g. keeping trailing zeros like
# "321.60" instead of Excel auto-converting to the number 321.6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.894.34'
$ws.Range("E2").Value = '  +0.73%  '

$ws.Range("D3").Value = '1.860.32'
$ws.Range("E3").Value = '  -0.10%  '

$ws.Range("E4").Value = '  -1.84%  '

$ws.Range("D5").Value = '''321.60'
$ws.Range("E5").Value = '  -0.84%  '

$ws.Range("E6").Value = '  -1.49%  '

$ws.Range("D7").Value = '''0.4334'
$ws.Range("E7").Value = '  -1.60%  '

$ws.Range("D8").Value = '''0.3800'
$ws.Range("E8").Value = '  +0.15%  '

$ws.Range("D9").Value = '''0.07431'
$ws.Range("E9").Value = '  -0.25%  '

$ws.Range("D10").Value = '''0.8870'
$ws.Range("E10").Value = '  +0.48%  '

$ws.Range("D11").Value = '''21.74'
$ws.Range("E11").Value = '  -0.12%  '

$ws.Range("D12").Value = '1.872.03'
$ws.Range("E12").Value = '  +0.17%  '

$ws.Range("D13").Value = '''6.779'
$ws.Range("E13").Value = '  +0.54%  '

$ws.Range("D14").Value = '''5.509'
$ws.Range("E14").Value = '  -0.68%  '

$ws.Range("D15").Value = '''0.07153'
$ws.Range("E15").Value = '  -0.70%  '

$ws.Range("D16").Value = '''88.48'
$ws.Range("E16").Value = '  +5.70%  '

$ws.Range("D17").Value = '''1.020'
$ws.Range("E17").Value = '  -1.59%  '

$ws.Range("D18").Value = '''0.000009059'
$ws.Range("E18").Value = '  -0.49%  '

$ws.Range("E19").Value = '  -1.55%  '

$ws.Range("D20").Value = '''15.57'
$ws.Range("E20").Value = '  +1.04%  '

$ws.Range("D21").Value = '27.926.79'
$ws.Range("E21").Value = '  +0.72%  '

$ws.Range("D22").Value = '''5.284'
$ws.Range("E22").Value = '  -0.48%  '

$ws.Range("D23").Value = '''11.23'
$ws.Range("E23").Value = '  -1.81%  '

$ws.Range("D24").Value = '2.094.53'
$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("D25").Value = '''2.032'
$ws.Range("E25").Value = '  +4.58%  '

$ws.Range("D26").Value = '''157.03'
$ws.Range("E26").Value = '  -0.46%  '

$ws.Range("D27").Value = '''18.72'
$ws.Range("E27").Value = '  -0.63%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '''2.028'
$ws.Range("E28").Value = '  +1.62%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '''5.440'
$ws.Range("E29").Value = '  +2.58%  '

$ws.Range("D30").Value = '''120.86'
$ws.Range("E30").Value = '  +2.97%  '

$ws.Range("D31").Value = '''0.08998'
$ws.Range("E31").Value = '  -0.62%  '

$ws.Range("D32").Value = '''1.242'
$ws.Range("E32").Value = '  +2.45%  '

$ws.Range("D33").Value = '''0.7779'
$ws.Range("E33").Value = '  +1.18%  '

$ws.Range("D34").Value = '''4.601'
$ws.Range("E34").Value = '  +0.60%  '

$ws.Range("D35").Value = '''2.925'
$ws.Range("E35").Value = '  -2.92%  '

$ws.Range("D36").Value = '''1.153'
$ws.Range("E36").Value = '  -0.79%  '

$ws.Range("D37").Value = '''1.017'
$ws.Range("E37").Value = '  -1.45%  '

$ws.Range("D38").Value = '''0.01977'
$ws.Range("E38").Value = '  -0.39%  '

$ws.Range("D39").Value = '''0.05330'
$ws.Range("E39").Value = '  -0.27%  '

$ws.Range("D40").Value = '''2.885'
$ws.Range("E40").Value = '  +1.95%  '

$ws.Range("D41").Value = '''0.5215'
$ws.Range("E41").Value = '  +0.46%  '

$ws.Range("D42").Value = '''7.028'
$ws.Range("E42").Value = '  +2.75%  '

$ws.Range("D43").Value = '''0.1685'
$ws.Range("E43").Value = '  -0.32%  '

$ws.Range("D44").Value = '''8.791'
$ws.Range("E44").Value = '  +1.12%  '

$ws.Range("D45").Value = '''111.05'
$ws.Range("E45").Value = '  +1.47%  '

$ws.Range("D46").Value = '''10.81'
$ws.Range("E46").Value = '  +1.64%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '''1.720'
$ws.Range("E47").Value = '  -0.36%  '

$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").Value = '''0.4765'
$ws.Range("E48").Value = '  +1.57%  '

$ws.Range("D49").Value = '''0.06497'
$ws.Range("E49").Value = '  +1.12%  '

$ws.Range("E50").Value = '  -1.43%  '

$ws.Range("D51").Value = '''1.885'
$ws.Range("E51").Value = '  +0.74%  '
